$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ----- Header row (row 1) -----
$ws.Range("A1").Value = "S.No."
$ws.Range("B1").Value = "ScenarioName"
$ws.Range("C1").Value = "resource"
$ws.Range("D1").Value = "api_Request"
$ws.Range("E1").Value = "key"
$ws.Range("F1").Value = "content_Type"
$ws.Range("G1").Value = "accuracy"
$ws.Range("H1").Value = "name"
$ws.Range("I1").Value = "phone_number"
$ws.Range("J1").Value = "address"
$ws.Range("K1").Value = "website"
$ws.Range("L1").Value = "language"
$ws.Range("M1").Value = "lat"
$ws.Range("N1").Value = "lng"
$ws.Range("O1").Value = "types"

# ----- Row 2 (non quote-prefixed cells first) -----
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Verify if place is being added using Add Place API"
$ws.Range("C2").Value = "/maps/api/place/add/json"
$ws.Range("D2").Value = "POST"
$ws.Range("E2").Value = "qaclick123"
$ws.Range("F2").Value = "JSON"
$ws.Range("H2").Value = "Frontline house"
$ws.Range("J2").Value = "India"
$ws.Range("K2").Value = "pranjal.com"
$ws.Range("L2").Value = "Hindi"
$ws.Range("O2").Value = "Shop#Shoe Park"

# ----- Row 3 -----
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Verify if place is not being added using Add Place API"
$ws.Range("C3").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = "qaclick123"
$ws.Range("F3").Value = "JSON"
$ws.Range("G3").Value = "NA"
$ws.Range("H3").Value = "NA"
$ws.Range("I3").Value = "NA"
$ws.Range("J3").Value = "NA"
$ws.Range("K3").Value = "NA"
$ws.Range("L3").Value = "NA"
$ws.Range("M3").Value = "NA"
$ws.Range("N3").Value = "NA"
$ws.Range("O3").Value = "NA"

# ----- Column widths -----
$ws.Columns.Item(2).ColumnWidth = 49
$ws.Columns.Item(3).ColumnWidth = 24.85546875
$ws.Columns.Item(4).ColumnWidth = 24.85546875
$ws.Columns.Item(5).ColumnWidth = 10
$ws.Columns.Item(6).ColumnWidth = 12.140625
$ws.Columns.Item(9).ColumnWidth = 11

# ----- Centered alignment style for column A (rows 2-17), applied once then
# copy/pasted down so the style engine doesn't leave orphan intermediate
# styles behind for each property assignment. Must run BEFORE the
# quote-prefixed text cells below so the resulting style indices match
# the target (alignment style = index 1, quote-prefix style = index 2). -----
$src = $ws.Range("A2")
$src.HorizontalAlignment = -4108
$src.VerticalAlignment = -4108
$src.Copy()
$dest = $ws.Range("A3:A17")
$dest.PasteSpecial(-4122)

# ----- Row 2 quote-prefixed (numeric-looking text) cells -----
$ws.Range("G2").Value = "'50"
$ws.Range("I2").Value = "'9822789334"
$ws.Range("M2").Value = "'-38.383494"
$ws.Range("N2").Value = "'33.427362"
